$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 84; this shifts rows 84..205 down to 85..206
$ws.Rows.Item(84).Insert(-4121)

# Populate the newly inserted row 84 with the new weekly record.
# Columns A,B,C,E,F,G,H,N,Q,R keep the same constant values used throughout
# the sheet for this market/product; D,I,J,K,L,M,O,P hold the new data point.
$ws.Cells.Item(84, 1).Value = 1
$ws.Cells.Item(84, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(84, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(84, 4).Value = 44495
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 5).Value = 15
$ws.Cells.Item(84, 6).Value = 100114013
$ws.Cells.Item(84, 7).Value = "Zanahoria"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 70
$ws.Cells.Item(84, 11).Value = 11000
$ws.Cells.Item(84, 12).Value = 12000
$ws.Cells.Item(84, 13).Value = 11500
$ws.Cells.Item(84, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(84, 15).Value = "Valle de Camiña"
$ws.Cells.Item(84, 16).Value = 460
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
